# Scheduled-runner refresh of cached market-board prices / leve profit
# calculations across the Bahamut_Profits workbook (ALC/ARM/BSM/CRP/CUL/
# GSM/LTW/WVR sheets). Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ)
# and LeveProfit(NQ/HQ) columns (H:N) for the rows whose market data moved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 1666.6666
$ws.Range("K31").Value = 4999.9998
$ws.Range("M31").Value = -4769.9998

$ws.Range("H129").Value = 962077.6
$ws.Range("I129").Value = 262
$ws.Range("J129").Value = 1393236.4
$ws.Range("K129").Value = 786
$ws.Range("L129").Value = 4179709.2
$ws.Range("M129").Value = 4214
$ws.Range("N129").Value = -4189709.2

$ws.Range("H132").Value = 2017.4517
$ws.Range("I132").Value = 2094.3103
$ws.Range("J132").Value = 903
$ws.Range("K132").Value = 6282.9309
$ws.Range("L132").Value = 2709
$ws.Range("M132").Value = -3752.9309
$ws.Range("N132").Value = -7769

$ws.Range("H135").Value = 2291.1765
$ws.Range("I135").Value = 1766.3077
$ws.Range("J135").Value = 3997
$ws.Range("K135").Value = 15896.7693
$ws.Range("L135").Value = 35973
$ws.Range("M135").Value = -13361.7693
$ws.Range("N135").Value = -41043

$ws.Range("H137").Value = 1000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7903.34
$ws.Range("I32").Value = 6594.839
$ws.Range("K32").Value = 6594.839
$ws.Range("M32").Value = -6307.839

$ws.Range("H61").Value = 1260.0646
$ws.Range("I61").Value = 1343.037
$ws.Range("K61").Value = 1343.037
$ws.Range("M61").Value = -1131.037

$ws.Range("H74").Value = 1354.8276
$ws.Range("I74").Value = 1349.6428
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1349.6428
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -475.6428000000001
$ws.Range("N74").Value = -3248

$ws.Range("H77").Value = 1354.8276
$ws.Range("I77").Value = 1349.6428
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 6748.214
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -2380.214
$ws.Range("N77").Value = -16236

$ws.Range("H122").Value = 1171.8667
$ws.Range("I122").Value = 1025
$ws.Range("J122").Value = 1269.7778
$ws.Range("K122").Value = 3075
$ws.Range("L122").Value = 3809.3334
$ws.Range("M122").Value = -625
$ws.Range("N122").Value = -8709.3334

$ws.Range("H132").Value = 1775.6852
$ws.Range("I132").Value = 1441.9524
$ws.Range("J132").Value = 2943.75
$ws.Range("K132").Value = 4325.857199999999
$ws.Range("L132").Value = 8831.25
$ws.Range("M132").Value = -1795.857199999999
$ws.Range("N132").Value = -13891.25

$ws.Range("H136").Value = 1260.0646
$ws.Range("I136").Value = 1343.037
$ws.Range("K136").Value = 4029.111
$ws.Range("M136").Value = -1479.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19117.271
$ws.Range("I134").Value = 1935.7317
$ws.Range("J134").Value = 58253
$ws.Range("K134").Value = 5807.1951
$ws.Range("L134").Value = 174759
$ws.Range("M134").Value = -3272.1951
$ws.Range("N134").Value = -179829

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2508.5818
$ws.Range("I31").Value = 2243.6924
$ws.Range("J31").Value = 3154.25
$ws.Range("K31").Value = 2243.6924
$ws.Range("L31").Value = 3154.25
$ws.Range("M31").Value = -1948.6924
$ws.Range("N31").Value = -3744.25

$ws.Range("H34").Value = 2508.5818
$ws.Range("I34").Value = 2243.6924
$ws.Range("J34").Value = 3154.25
$ws.Range("K34").Value = 2243.6924
$ws.Range("L34").Value = 3154.25
$ws.Range("M34").Value = -2041.6924
$ws.Range("N34").Value = -3558.25

$ws.Range("H58").Value = 3188.2766
$ws.Range("I58").Value = 875.6
$ws.Range("K58").Value = 875.6
$ws.Range("M58").Value = -672.6

$ws.Range("H136").Value = 3188.2766
$ws.Range("I136").Value = 875.6
$ws.Range("K136").Value = 2626.8
$ws.Range("M136").Value = -76.80000000000018

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 200000080
$ws.Range("I40").Value = 130
$ws.Range("K40").Value = 520
$ws.Range("M40").Value = -451

$ws.Range("H97").Value = 5952706.5
$ws.Range("I97").Value = 7143193
$ws.Range("J97").Value = 275
$ws.Range("K97").Value = 21429579
$ws.Range("L97").Value = 825
$ws.Range("M97").Value = -21429083
$ws.Range("N97").Value = -1817

$ws.Range("H98").Value = 3619.7778
$ws.Range("I98").Value = 600
$ws.Range("J98").Value = 3997.25
$ws.Range("K98").Value = 1800
$ws.Range("L98").Value = 11991.75
$ws.Range("M98").Value = -302
$ws.Range("N98").Value = -14987.75

$ws.Range("H107").Value = 433094.72
$ws.Range("I107").Value = 1585
$ws.Range("J107").Value = 648849.5600000001
$ws.Range("K107").Value = 4755
$ws.Range("L107").Value = 1946548.68
$ws.Range("M107").Value = -2835
$ws.Range("N107").Value = -1950388.68

$ws.Range("H132").Value = 872.2222
$ws.Range("I132").Value = 846.6667
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 7620.0003
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5090.0003
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3924.7917
$ws.Range("I80").Value = 3849.8572
$ws.Range("J80").Value = 4029.7
$ws.Range("K80").Value = 3849.8572
$ws.Range("L80").Value = 4029.7
$ws.Range("M80").Value = -2851.8572
$ws.Range("N80").Value = -6025.7

$ws.Range("H83").Value = 3924.7917
$ws.Range("I83").Value = 3849.8572
$ws.Range("J83").Value = 4029.7
$ws.Range("K83").Value = 19249.286
$ws.Range("L83").Value = 20148.5
$ws.Range("M83").Value = -14257.286
$ws.Range("N83").Value = -30132.5

$ws.Range("H97").Value = 1270.0322
$ws.Range("I97").Value = 1393.6364
$ws.Range("K97").Value = 1393.6364
$ws.Range("M97").Value = -897.6364000000001

$ws.Range("H113").Value = 4030.3684
$ws.Range("I113").Value = 4420.1665
$ws.Range("J113").Value = 3362.1428
$ws.Range("K113").Value = 4420.1665
$ws.Range("L113").Value = 3362.1428
$ws.Range("M113").Value = -2250.1665
$ws.Range("N113").Value = -7702.1428

$ws.Range("H136").Value = 8474.695
$ws.Range("J136").Value = 8474.695
$ws.Range("L136").Value = 25424.085
$ws.Range("N136").Value = -30524.085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 760.94116
$ws.Range("I113").Value = 697.5
$ws.Range("J113").Value = 913.2
$ws.Range("K113").Value = 2092.5
$ws.Range("L113").Value = 2739.6
$ws.Range("M113").Value = 77.5
$ws.Range("N113").Value = -7079.6
